$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Three new observers were appended to the bottom of the observer_ids table
# (row 29's "observer_id" cell was typed as the text "ZS" instead of a
# number in the source data, so it is written as text here too).
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "OM"
$ws.Range("C27").Value = 2

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "NI"
$ws.Range("C28").Value = 2

$ws.Range("A29").Value = "ZS"
$ws.Range("B29").Value = "ZS"
$ws.Range("C29").Value = 1

# Match the explicit (customHeight) row height already used by every
# other data row in the sheet.
$ws.Rows.Item(27).RowHeight = 23.25
$ws.Rows.Item(28).RowHeight = 23.25
$ws.Rows.Item(29).RowHeight = 23.25

# Reproduce the author's final selection state as closely as possible:
# click A7, then select the whole sheet (Ctrl+A).
$ws.Range("A7").Select()
$ws.Cells.Select()
